$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: updated publish date
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# (Setting a plain "true" string via .Value auto-coerces to a Boolean in
#  Excel, so stage the literal text through a formula cell and paste the
#  computed value back in, which keeps it as real text.)
$ws.Range("D1").Formula = '="true"'
$ws.Range("D1").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4163) | Out-Null
$ws.Range("D1").Clear() | Out-Null
